# Update gh-pages to output generated at 456a3b4
#
# The upstream scraper re-ran and produced refreshed "want to go" counts
# (column F) for several existing events, and discovered one brand-new
# event that needs to be appended as a new row, on both the "展览" and
# "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Updated "想去人数" (want-to-go count) values -----------------
    $ws.Range("F2").Value  = 1070
    $ws.Range("F5").Value  = 3068
    $ws.Range("F7").Value  = 2305
    $ws.Range("F11").Value = 1109
    $ws.Range("F15").Value = 540
    $ws.Range("F19").Value = 10

    # --- New row 23: 九江·第三届ACD动漫游戏嘉年华 -----------------------
    # Duplicate row 22 first so the new row inherits identical cell
    # formatting/styles (e.g. the bordered/bold style used in column A),
    # then overwrite the cell contents with the new record's data.
    $ws.Range("A22:I22").Copy($ws.Range("A23:I23"))

    $ws.Range("A23").Value = 22

    # "2024-05-01" looks like a date, so force the cell to Text first to
    # keep it stored as the literal string (matching how the rest of the
    # sheet stores dates as plain text), then restore a plain/default
    # format so the cell's visual style matches its neighbours.
    $ws.Range("B23").NumberFormat = "@"
    $ws.Range("B23").Value = "2024-05-01"
    $ws.Range("B23").Style = "Normal"

    $ws.Range("C23").Value = "九江·第三届ACD动漫游戏嘉年华"
    $ws.Range("D23").Value = "九瑞大道与重庆路交汇处西南角 九江国际会展中心"
    $ws.Range("E23").Value = "2024.05.01 09:00-05.02 17:00"
    $ws.Range("F23").Value = 1
    $ws.Range("G23").Value = 39.9
    $ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=82464"
    $ws.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202403/HjMMyP3a1709780146797.jpeg"
}
